# Add two new columns, I (I0) and J (IF), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I1 and J1, formatted like the existing header H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Data values for the new columns.
$values = @{
    2  = @(9, 9)
    3  = @(10, 10)
    4  = @(6, 6)
    5  = @(7, 8)
    6  = @(7, 7)
    7  = @(7, 8)
    8  = @(7, 8)
    9  = @(6, 7)
    10 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
